$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset transaction counts for A2:G2 to 1000 (H2 and I2 already 1000)
$ws.Range("A2:G2").Value = 1000
